# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (column G) values, recalculated (replacing the old Strike# based values)
$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 2
    10 = 0
    11 = 0
    12 = 1
    13 = 2
    14 = 0
    15 = 2
    16 = 0
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 2
    25 = 0
    26 = 1
    27 = 2
    28 = 1
    29 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
